$wb = $excel.ActiveWorkbook

# --- ACHData (TypeOfAccount data): row 3 (ID=2) AccountNumber / ConfirmAccountNumber
# both corrected to the shared value "95125489" (previously mismatched
# "65987541" / "65239412", which become orphaned and are dropped from the
# shared-string table on save).
$wsACH = $wb.Worksheets.Item("ACHData")
$wsACH.Range("B3").Value = "95125489"
$wsACH.Range("C3").Value = "95125489"

# --- UDFData: add the BWP Test Cases and Keywords rows (ID=2 and ID=3)
$wsUDF = $wb.Worksheets.Item("UDFData")

$wsUDF.Range("B3").Value = "udf data 1"
$wsUDF.Range("B3").Style = "Normal"
$wsUDF.Range("C3").Value = "udf data 2"
$wsUDF.Range("C3").Style = "Normal"
$wsUDF.Range("D3").Value = "udf data 3"
$wsUDF.Range("D3").Style = "Normal"
$wsUDF.Range("E3").Value = "Sweet"
$wsUDF.Range("E3").Style = "Normal"
$wsUDF.Range("F3").Value = "Sour"
$wsUDF.Range("F3").Style = "Normal"
$wsUDF.Range("G3").Value = "udf data 6"
$wsUDF.Range("G3").Style = "Normal"
$wsUDF.Range("H3").Value = "udf data 7"
$wsUDF.Range("H3").Style = "Normal"
$wsUDF.Range("I3").Value = "udf data 8"
$wsUDF.Range("I3").Style = "Normal"
$wsUDF.Range("J3").Value = "udf data 9"
$wsUDF.Range("J3").Style = "Normal"
$wsUDF.Range("K3").Value = "udf data 10"
$wsUDF.Range("K3").Style = "Normal"

$wsUDF.Range("A4").Value = "3"

# --- Selection / active-sheet bookkeeping to mirror the saved UI state ---

# NameData: cursor moved from F1 to H11
$wsName = $wb.Worksheets.Item("NameData")
$wsName.Range("H11").Select()

# UDFData: cursor moved from J2 to A4 (no longer the active tab)
$wsUDF.Range("A4").Select()

# ACHData becomes the active tab, cursor on C3
$wsACH.Activate()
$wsACH.Range("C3").Select()
